$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows describing the new user story / task
$ws.Range("A21").Value = "Sistem treba da limitira maksimalni iznos uplate/isplate po novcaniku na mesecnom nivou"

$ws.Range("B22").Value = "Dodati max/min ogranicenja u appsettings.json kao njihovu proveru u PayIn/PayOut metode"
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 10

# Update the selection to mirror the author's final cursor position
$ws.Range("C27").Select()
